$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current row 395, shifting rows 395-414 down to 400-419.
$ws.Rows("395:399").Insert()

# Common (unchanged) columns for this block of records.
$marketId   = 6
$marketName = "Mercado Mayorista Lo Valledor de Santiago"
$region     = "Metropolitana"
$codreg     = 13
$tipo       = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria  = "Ciruela"

# New row data: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Origen, PrecioKg, KgUnidad
$rows = @(
    @{ Row=395; Fecha=44585; Variedad="Black Amber"; Calidad="Especial"; Volumen=170; PMin=13000;  PMax=13000;  PProm=13000;  Unidad="$/caja 15 kilos granel"; Origen="Región de O'Higgins"; PrecioKg=867; KgUnidad=15 },
    @{ Row=396; Fecha=44585; Variedad="Black Amber"; Calidad="Primera";  Volumen=15;  PMin=210000; PMax=210000; PProm=210000; Unidad="$/bins (450 kilos)";       Origen="Región de O'Higgins"; PrecioKg=467; KgUnidad=450 },
    @{ Row=397; Fecha=44585; Variedad="Black Amber"; Calidad="Primera";  Volumen=300; PMin=10000;  PMax=11000;  PProm=10500;  Unidad="$/caja 15 kilos granel"; Origen="Región de O'Higgins"; PrecioKg=700; KgUnidad=15 },
    @{ Row=398; Fecha=44585; Variedad="Black Amber"; Calidad="Segunda";  Volumen=12;  PMin=160000; PMax=160000; PProm=160000; Unidad="$/bins (450 kilos)";       Origen="Región de O'Higgins"; PrecioKg=356; KgUnidad=450 },
    @{ Row=399; Fecha=44585; Variedad="Black Amber"; Calidad="Segunda";  Volumen=200; PMin=8000;   PMax=8000;   PProm=8000;   Unidad="$/caja 15 kilos granel"; Origen="Región de O'Higgins"; PrecioKg=533; KgUnidad=15 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value  = $marketId
    $ws.Cells.Item($i, 2).Value  = $marketName
    $ws.Cells.Item($i, 3).Value  = $region
    $ws.Cells.Item($i, 4).Value  = $r.Fecha
    $ws.Cells.Item($i, 5).Value  = $codreg
    $ws.Cells.Item($i, 6).Value  = $tipo
    $ws.Cells.Item($i, 7).Value  = $productoId
    $ws.Cells.Item($i, 8).Value  = $producto
    $ws.Cells.Item($i, 9).Value  = $categoriaId
    $ws.Cells.Item($i, 10).Value = $categoria
    $ws.Cells.Item($i, 11).Value = $r.Variedad
    $ws.Cells.Item($i, 12).Value = $r.Calidad
    $ws.Cells.Item($i, 13).Value = $r.Volumen
    $ws.Cells.Item($i, 14).Value = $r.PMin
    $ws.Cells.Item($i, 15).Value = $r.PMax
    $ws.Cells.Item($i, 16).Value = $r.PProm
    $ws.Cells.Item($i, 17).Value = $r.Unidad
    $ws.Cells.Item($i, 18).Value = $r.Origen
    $ws.Cells.Item($i, 19).Value = $r.PrecioKg
    $ws.Cells.Item($i, 20).Value = $r.KgUnidad
}
